$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell F2 (rpc-reply message-id uuid changed)
$ws.Range("F2").Value = @"
<rpc-reply message-id="urn:uuid:9552cef2-111e-4c24-b16f-c81e38dc1766">
  <data/>
</rpc-reply>

"@

# Update cell G2 (identifier/name changed from BGP_65000 protocol entries to default with oc-pol-types prefix)
$ws.Range("G2").Value = @"
<edit-config>
    <target>
      <candidate/>
    </target>
    <config>
      <network-instances xmlns="http://openconfig.net/yang/network-instance">
        <network-instance>
          <name>Prueba_LxVPN</name>
          <config>
            <name>Prueba_LxVPN</name>
            <type xmlns:oc-ni-types="http://openconfig.net/yang/network-instance-types">oc-ni-types:L3VRF</type>
          </config>
          <protocols>
            <protocol>
              <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
              <name>default</name>
              <config>
                <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
                <name>default</name>
              </config>
              <bgp>
                <global>
                  <config>
                    <as>65000</as>
                  </config>
                </global>
                <neighbors>
                  <neighbor>
                    <neighbor-address>192.168.1.2</neighbor-address>
                    <config>
                      <neighbor-address>192.168.1.2</neighbor-address>
                      <peer-type>EXTERNAL</peer-type>
                    </config>
                  </neighbor>
                </neighbors>
              </bgp>
            </protocol>
          </protocols>
        </network-instance>
      </network-instances>
    </config>
  </edit-config>
"@
